# [fixed_asset template] missing subtotal
#
# Adds a per-group subtotal row (F/G/H = total book value from / accum dep /
# total book value to) just before the "{{/each}}" that closes the groups
# loop, and makes the grand-total row's totals + "Total" label bold to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank row at row 12 -----------------------------------
# This pushes the old row 12 ("{{/each}}" closing the groups loop) down to
# row 13, and the old row 13 (grand Total row) down to row 14, updating the
# merged-cell reference (A13:E13 -> A14:E14) automatically.
$ws.Rows.Item(12).Insert()

# --- 2. Build the bold "totals" look (no fill/border, just bold Arial 10) --
# Grab an existing plain cell's format, then flip Bold on - this reuses the
# workbook's existing font/fill table instead of inventing new ones.
$ws.Range("F13").Copy()
$ws.Range("F12:H12").PasteSpecial(-4122)
$ws.Range("F12:H12").Font.Bold = $true

# --- 3. Fill in the new per-group subtotal row (row 12) --------------------
$ws.Range("F12").Value = "{{currency total_book_val_from}}"
$ws.Range("G12").Value = "{{currency total_accum_dep}}"
$ws.Range("H12").Value = "{{currency total_book_val_to}}"

# --- 4. Make the grand-total row's numbers bold too (row 14, was row 13) --
$ws.Range("F14:H14").Font.Bold = $true

# --- 5. Make the grand-total "Total" label bold (A14:E14, was A13:E13) ----
$ws.Range("A14:E14").Font.Bold = $true

$wb.Save()
